$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember the hyperlinked cells (column G/I, rows 2-8) and the text they
# display *before* the shift - that text is itself the link target, so it
# lets us recreate the hyperlinks after everything moves one column right.
$hyperlinkCols = @("G", "I")
$targets = @{}
foreach ($col in $hyperlinkCols) {
    for ($r = 2; $r -le 8; $r++) {
        $addr = "$col$r"
        $targets[$addr] = $ws.Range($addr).Text
    }
}

# Insert a new blank column before column A, shifting all existing
# columns (and the header row) one place to the right.
$ws.Range("A:A").Insert()

# The worksheet's Hyperlinks collection does not follow the shifted
# cells automatically, so rebuild it at the new (one-column-over)
# locations. Deleting through any single-cell filter clears the whole
# sheet collection here, which is fine since we are about to re-add
# every entry anyway.
$ws.Range("A1").Hyperlinks.Delete()

for ($r = 2; $r -le 8; $r++) {
    foreach ($col in $hyperlinkCols) {
        $newCol = [char]([int][char]$col + 1)
        $oldAddr = "$col$r"
        $newAddr = "$newCol$r"
        $ws.Hyperlinks.Add($ws.Range($newAddr), $targets[$oldAddr]) | Out-Null
    }
}

# Fill the new column A (rows 2-8) with a zero-based row index, copying
# the header's style (bold/centered/bordered) onto those cells.
$ws.Range("B1").Copy()
$ws.Range("A2:A8").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

for ($i = 0; $i -le 6; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $i
}
